# AVL tree time experiment - final changes to new fork
#
# Adds a new "Unnamed: 0" column (E) that mirrors the original row index,
# and refreshes the measured timing values for the Start / Random / End
# columns (which are now columns C, D, E after the new column is inserted
# before them, at B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell E1 ("Unnamed: 0" header slot shifts the existing
#     Start/Random/End headers over - duplicate D1's look/format onto E1) ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E1").Value = "End"

# B1/C1/D1 keep their original shared-string slots; only the E1 label is new.
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Start"
$ws.Range("D1").Value = "Random"

# --- Data rows: column B becomes the plain row index (0..9, same as A);
#     columns C, D, E get the refreshed timing measurements ---
$rows = @(
    @{ Row = 2;  B = 0; C = "7.572733333320986E-05";  D = "8.261439999993554E-05";  E = "8.70215999999952E-05"  },
    @{ Row = 3;  B = 1; C = "7.810396666665533E-05";  D = "8.707593333330503E-05";  E = "8.787199999998542E-05" },
    @{ Row = 4;  B = 2; C = "8.046437777774272E-05";  D = "8.749302222223478E-05";  E = "9.223591111109879E-05" },
    @{ Row = 5;  B = 3; C = "8.39691666666719E-05";   D = "8.9694799999999E-05";    E = "9.211456666666133E-05" },
    @{ Row = 6;  B = 4; C = "8.478923999997884E-05";  D = "9.173118666667506E-05";  E = "9.651869333335222E-05" },
    @{ Row = 7;  B = 5; C = "8.435133333334004E-05";  D = "9.18643777777864E-05";   E = "9.165940000000066E-05" },
    @{ Row = 8;  B = 6; C = "8.51497809523954E-05";   D = "9.525457142858401E-05";  E = "9.690323809523506E-05" },
    @{ Row = 9;  B = 7; C = "8.641094166665651E-05";  D = "9.751739166665634E-05";  E = "0.0001019106416666773" },
    @{ Row = 10; B = 8; C = "9.483948888888942E-05";  D = "0.0001178310740740724";  E = "0.0001164564074074069" },
    @{ Row = 11; B = 9; C = "0.0001009554333333351";  D = "0.0001033752199999981";  E = "0.0001082831466666751" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = [double]$r.C
    $ws.Cells.Item($r.Row, 4).Value = [double]$r.D
    $ws.Cells.Item($r.Row, 5).Value = [double]$r.E
}
